$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6614
$ws.Range("B2").Value = 45946.40625

$ws.Range("A3").Value = 6483
$ws.Range("B3").Value = 45946.41666666666

$ws.Range("A4").Value = 6386
$ws.Range("B4").Value = 45946.42708333334

$ws.Range("B2:B4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
